$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-13 Saturday", "2024-07-14 Sunday"),
    @("998÷2=", "928÷6="),
    @("802÷6=", "458÷5="),
    @("367÷2=", "560÷9="),
    @("740÷3=", "307÷4="),
    @("924÷3=", "360÷6="),
    @("115÷3=", "166÷4="),
    @("643÷6=", "730÷7="),
    @("179÷3=", "108÷2="),
    @("452÷8=", "919÷8="),
    @("527÷7=", "757÷2="),
    @("145÷7=", "163÷3="),
    @("858÷3=", "683÷2="),
    @("176÷4=", "873÷8="),
    @("767÷3=", "816÷6="),
    @("174÷2=", "324÷9="),
    @("636÷2=", "486÷3="),
    @("965÷2=", "208÷7="),
    @("176÷9=", "122÷7="),
    @("722÷9=", "116÷2="),
    @("664÷2=", "345÷3="),
    @("772÷3=", "349÷7="),
    @("241÷6=", "397÷8="),
    @("686÷2=", "713÷6="),
    @("606÷9=", "770÷4="),
    @("231÷6=", "252÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
